$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new day's mod-count entry as row 17.
# The date is entered with a leading apostrophe so Excel keeps it as the
# literal text "2025/11/26" instead of auto-converting it to a date serial.
$ws.Range("A17").Value = "'2025/11/26"
$ws.Range("B17").Value = "逃离鸭科夫"
$ws.Range("C17").Value = 1263

# Match the formatting of the preceding data rows (centered alignment style).
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
